$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, pushing the existing data rows down by one
# (new IPO entry added above the rest of the RPA dataset rows).
$ws.Rows.Item(2).Insert()

# The "date" columns (A, O, P) are stored as plain text in this sheet, not
# as real Excel dates, so prefix with an apostrophe to force text and avoid
# auto-conversion to a date serial number.
$ws.Range("A2").Value = "'2023-11-07"
$ws.Range("B2").Value = "쏘닉스"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 270
$ws.Range("E2").Value = "KB"
$ws.Range("F2").Value = 270
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 7500
$ws.Range("N2").Value = 100
$ws.Range("O2").Value = "'2023-10-26"
$ws.Range("P2").Value = "'2023-10-31"
$ws.Range("Q2").Value = 2600000

# Excel's Insert() copies the formatting of the row above (the bold header
# row here). The other data rows carry no explicit style, so strip that
# back off the freshly inserted row to match the rest of the sheet.
$ws.Rows.Item(2).ClearFormats()
